$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (entry "A 27984-2024"); remaining rows shift up by one.
$ws.Rows.Item(2).Delete()

# After the shift, the table now occupies rows 2..31 (was 2..32).
# Delete the trailing 3 rows that are no longer present (previously rows 30,31,32;
# now, after the first delete, they are rows 29,30,31).
$ws.Range("A29:Z31").EntireRow.Delete()

# Update column C ("Förändrad") to 45478 for the remaining data rows (2..28).
$ws.Range("C2:C28").Value = 45478

# The new last row (28) should lose its explicit row height / customHeight flag,
# matching the original last row's (32) auto-height state.
$ws.Rows.Item(28).AutoFit()
